# Calculator.xlsx - "Worked on more Statistics"
# Adds a t-Statistics section (paired t-test / t-score / Cohen's d, plus
# independent-samples pooled-variance t-stat) below the existing Z-score
# calculator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Section header -------------------------------------------------
$ws.Range("A12").Value = "t Statistics"

# ---- Population / Sample details sub-headers -------------------------
$ws.Range("A13").Value = "Population Details"
$ws.Range("A13:B13").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A13:B13").Merge()

$ws.Range("D13").Value = "Sample Population Details"
$ws.Range("D13:E13").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D13:E13").WrapText = $true
$ws.Range("D13:E13").Merge()

# ---- Population / sample stats ---------------------------------------
$ws.Range("A14").Value = "Mean"
$ws.Range("B14").Value = 1830

$ws.Range("D14").Value = "Sample Mean"
$ws.Range("E14").Value = 1700

$ws.Range("D15").Value = "Std Dev"
$ws.Range("E15").Value = 200

$ws.Range("D16").Value = "Sample Size"
$ws.Range("E16").Value = 25

# ---- t Score / Cohen's d ---------------------------------------------
$ws.Range("B18").Value = "t Score"
$ws.Range("C18").Formula = "=(E14-B14)/(E15/SQRT(E16))"

$ws.Range("B19").Value = "Cohen's d"
$ws.Range("C19").Formula = "=(E14-B14)/E15"

# ---- Two-sample (pooled variance) t-test ------------------------------
# Labels are entered column-by-column (matches the original authoring
# order captured by the shared-strings table), values filled in after.
$ws.Range("A24").Value = "X1"
$ws.Range("A25").Value = "X2"

$ws.Range("D24").Value = "n1"
$ws.Range("D25").Value = "n2"

$ws.Range("G24").Value = "SS1"
$ws.Range("G25").Value = "SS2"

$ws.Range("B24").Value = 35.8
$ws.Range("B25").Value = 31.6

$ws.Range("E24").Value = 207
$ws.Range("E25").Value = 220

$ws.Range("H24").Value = 481
$ws.Range("H25").Value = 322

$ws.Range("A27").Value = "Sp^2"
$ws.Range("A29").Value = "t stats"
$ws.Range("D27").Value = "df"

$ws.Range("B27").Formula = "=(H24+H25) / (E27)"
$ws.Range("E27").Formula = "=E24+E25-2"

$ws.Range("A28").Value = "Std Error"
$ws.Range("B28").Formula = "=SQRT( ( (B27^2) / E24) + ( (B27^2) / E25) )"

$ws.Range("B29").Formula = "=(B24-B25)/B28"

# ---- Selection matches the author's final cursor position -------------
$ws.Range("B28").Select()
